# Add the "ObjTables" entry to the curated standards table, and update the
# sheet/window view state (zoom, split position, selection) to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append the new data row (row 46) -------------------------------------
$ws.Range("A46").Value = "ObjTables"
$ws.Range("B46").Value = "Biochemical data format"
$ws.Range("C46").Value = "Structured spreadsheets with ObjTables enable data reuse and integration"

# --- View-state tweaks (zoom + split/selection) ----------------------------
$win = $excel.ActiveWindow
$win.Zoom = 130
$win.SplitRow = 27

[void]$ws.Range("C48").Select()
